$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Количество предложений в тексте (ед.)"
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = "Cell is empty"

$ws.Range("A3").Value = "Количество слов в тексте (ед.)"
$ws.Range("B3").Value = 777
$ws.Range("C3").Value = "Cell is empty"

$ws.Range("A4").Value = "Абсолютное количество знаков препинаний в тексте (ед.)"
$ws.Range("B4").Value = 136
$ws.Range("C4").Value = "Cell is empty"

$ws.Range("A5").Value = "Относительное количество знаков препинаний в тексте (%)"
$ws.Range("B5").Value = 17.5
$ws.Range("C5").Value = "Cell is empty"

$ws.Range("A6").Value = "Среднее число слов в предложении (ед.)"
$ws.Range("B6").Value = 32.38
$ws.Range("C6").Value = "Cell is empty"

$ws.Range("A7").Value = "Среднее число букв в слове (ед.)"
$ws.Range("B7").Value = 6.8
$ws.Range("C7").Value = "Cell is empty"
